$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.857.26'
$ws.Cells.Item(2, 5).Value = '  +0.25%  '
$ws.Cells.Item(3, 4).Value = '2.298.59'
$ws.Cells.Item(3, 5).Value = '  -0.76%  '
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '102.99'
$ws.Cells.Item(5, 5).Value = '  +5.74%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '270.21'
$ws.Cells.Item(6, 5).Value = '  -0.28%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.626'
$ws.Cells.Item(7, 5).Value = '  +0.09%  '
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.610'
$ws.Cells.Item(9, 5).Value = '  -1.95%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '45.57'
$ws.Cells.Item(10, 5).Value = '  -0.29%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0937'
$ws.Cells.Item(11, 5).Value = '  -1.00%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.98'
$ws.Cells.Item(12, 5).Value = '  -1.72%  '
$ws.Cells.Item(13, 5).Value = '  +1.77%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '15.81'
$ws.Cells.Item(14, 5).Value = '  +2.22%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.860'
$ws.Cells.Item(15, 5).Value = '  -0.86%  '
$ws.Cells.Item(16, 4).Value = '2.301.46'
$ws.Cells.Item(16, 5).Value = '  -1.08%  '
$ws.Cells.Item(17, 4).Value = '43.744.36'
$ws.Cells.Item(17, 5).Value = '  +0.06%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.0000111'
$ws.Cells.Item(18, 5).Value = '  +1.83%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.26'
$ws.Cells.Item(19, 5).Value = '  -2.40%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '72.34'
$ws.Cells.Item(20, 5).Value = '  -0.47%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '2.49'
$ws.Cells.Item(21, 5).Value = '  +9.16%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '233.61'
$ws.Cells.Item(22, 5).Value = '  -2.43%  '
$ws.Cells.Item(24, 5).Value = '  -1.85%  '
$ws.Cells.Item(25, 5).Value = '  +0.02%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '11.22'
$ws.Cells.Item(26, 5).Value = '  -0.67%  '
$ws.Cells.Item(27, 5).Value = '  -0.42%  '
$ws.Cells.Item(28, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '39.85'
$ws.Cells.Item(28, 5).Value = '  +4.48%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.29'
$ws.Cells.Item(29, 5).Value = '  +0.54%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '177.57'
$ws.Cells.Item(30, 5).Value = '  +1.39%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '21.84'
$ws.Cells.Item(31, 5).Value = '  -2.47%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0900'
$ws.Cells.Item(32, 5).Value = '  -0.07%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.46'
$ws.Cells.Item(33, 5).Value = '  -0.24%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.84'
$ws.Cells.Item(34, 5).Value = '  +10.24%  '
$ws.Cells.Item(35, 5).Value = '  +0.34%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.109'
$ws.Cells.Item(36, 5).Value = '  +0.57%  '
$ws.Cells.Item(37, 5).Value = '  -1.76%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.56'
$ws.Cells.Item(38, 5).Value = '  +5.97%  '
$ws.Cells.Item(39, 5).Value = '  -0.15%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.236'
$ws.Cells.Item(40, 5).Value = '  -2.95%  '
$ws.Cells.Item(41, 5).Value = '  +1.57%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '12.36'
$ws.Cells.Item(42, 5).Value = '  +1.75%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '65.17'
$ws.Cells.Item(43, 5).Value = '  +5.10%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '8.83'
$ws.Cells.Item(44, 5).Value = '  -3.75%  '
$ws.Cells.Item(45, 5).Value = '  -2.17%  '
$ws.Cells.Item(46, 5).Value = '  -0.95%  '
$ws.Cells.Item(47, 5).Value = '  +1.08%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '98.41'
$ws.Cells.Item(48, 5).Value = '  -1.82%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.450'
$ws.Cells.Item(49, 5).Value = '  +8.56%  '
$ws.Cells.Item(50, 5).Value = '  +12.25%  '
$ws.Cells.Item(51, 4).Value = '2.522.19'
$ws.Cells.Item(51, 5).Value = '  -0.87%  '
